$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab title)
$ws.Name = "Through 2022-03-10"

# Update the workbook title text used for the sheet entry / label cell
$ws.Range("A4").Value = "March (through 03-10)"

# Update March row (row 4) values
$ws.Range("C4").Value = 14
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 17
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 45

# Update Total row (row 5) values
$ws.Range("C5").Value = 101
$ws.Range("D5").Value = 148
$ws.Range("E5").Value = 154
$ws.Range("G5").Value = 161
$ws.Range("H5").Value = 372
$ws.Range("I5").Value = 346
